# Updated cryptos list (GitHub Actions scrape refresh): new Price/Volume(1h)
# values per coin, plus Chainlink/Polkadot swapping rank positions 14/15.
# Numeric-looking Price strings are written with a leading apostrophe so
# Excel keeps them as text (matching the original inlineStr cells) instead
# of auto-coercing to numbers; the style is then reset to "Normal" so no
# quotePrefix/number-format drift is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.558.29'
$ws.Range('E2').Value = '  +3.49%  '
$ws.Range('D3').Value = '1.795.34'
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('D4').Value = '''1.003'
$ws.Range('E4').Value = '  +0.30%  '
$ws.Range('D5').Value = '''313.54'
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('E6').Value = '  +0.21%  '
$ws.Range('D7').Value = '''0.5298'
$ws.Range('E7').Value = '  -0.76%  '
$ws.Range('D8').Value = '''0.3767'
$ws.Range('E8').Value = '  -0.27%  '
$ws.Range('D9').Value = '''0.07526'
$ws.Range('E9').Value = '  -0.20%  '
$ws.Range('D10').Value = '''42.58'
$ws.Range('E10').Value = '  -0.66%  '
$ws.Range('D11').Value = '''1.118'
$ws.Range('E11').Value = '  +0.25%  '
$ws.Range('D12').Value = '''21.14'
$ws.Range('E12').Value = '  +0.81%  '
$ws.Range('D13').Value = '''1.003'
$ws.Range('E13').Value = '  +0.33%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = '''7.513'
$ws.Range('E14').Value = '  +5.90%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '''6.198'
$ws.Range('E15').Value = '  +0.09%  '
$ws.Range('D16').Value = '1.797.44'
$ws.Range('E16').Value = '  +0.64%  '
$ws.Range('D17').Value = '''90.24'
$ws.Range('E17').Value = '  -0.82%  '
$ws.Range('D18').Value = '''0.00001067'
$ws.Range('E18').Value = '  -0.15%  '
$ws.Range('D19').Value = '''0.06463'
$ws.Range('E19').Value = '  -0.73%  '
$ws.Range('E20').Value = '  +0.19%  '
$ws.Range('D21').Value = '''17.29'
$ws.Range('E21').Value = '  +1.85%  '
$ws.Range('D22').Value = '''5.922'
$ws.Range('E22').Value = '  -0.40%  '
$ws.Range('D23').Value = '28.600.74'
$ws.Range('E23').Value = '  +3.52%  '
$ws.Range('D24').Value = '''11.19'
$ws.Range('E24').Value = '  -0.54%  '
$ws.Range('D25').Value = '''2.101'
$ws.Range('E25').Value = '  +0.41%  '
$ws.Range('D26').Value = '''160.89'
$ws.Range('E26').Value = '  +3.49%  '
$ws.Range('D27').Value = '''20.52'
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('D28').Value = '2.003.87'
$ws.Range('E28').Value = '  +0.48%  '
$ws.Range('D29').Value = '''2.359'
$ws.Range('E29').Value = '  -1.89%  '
$ws.Range('D30').Value = '''124.28'
$ws.Range('E30').Value = '  +1.73%  '
$ws.Range('D31').Value = '''1.117'
$ws.Range('E31').Value = '  -0.34%  '
$ws.Range('D32').Value = '''0.1031'
$ws.Range('E32').Value = '  +0.09%  '
$ws.Range('D33').Value = '''5.698'
$ws.Range('E33').Value = '  -0.11%  '
$ws.Range('D34').Value = '''3.677'
$ws.Range('E34').Value = '  +1.59%  '
$ws.Range('D35').Value = '''0.2285'
$ws.Range('E35').Value = '  +9.43%  '
$ws.Range('D36').Value = '''0.06543'
$ws.Range('E36').Value = '  +8.33%  '
$ws.Range('D37').Value = '''0.02319'
$ws.Range('E37').Value = '  +1.37%  '
$ws.Range('D38').Value = '''8.858'
$ws.Range('E38').Value = '  +2.52%  '
$ws.Range('D39').Value = '''5.068'
$ws.Range('E39').Value = '  +1.13%  '
$ws.Range('D40').Value = '''11.45'
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('D41').Value = '''0.6287'
$ws.Range('E41').Value = '  +0.43%  '
$ws.Range('D42').Value = '''1.213'
$ws.Range('E42').Value = '  +5.62%  '
$ws.Range('E43').Value = '  +0.18%  '
$ws.Range('E44').Value = '  -1.22%  '
$ws.Range('D45').Value = '''13.45'
$ws.Range('E45').Value = '  +0.17%  '
$ws.Range('D46').Value = '''0.5915'
$ws.Range('E46').Value = '  +0.52%  '
$ws.Range('D47').Value = '''3.670'
$ws.Range('E47').Value = '  +0.88%  '
$ws.Range('D48').Value = '''126.31'
$ws.Range('E48').Value = '  +3.52%  '
$ws.Range('D49').Value = '''1.976'
$ws.Range('E49').Value = '  +2.88%  '
$ws.Range('D50').Value = '''1.166'
$ws.Range('E50').Value = '  +2.87%  '
$ws.Range('D51').Value = '''0.06924'
$ws.Range('E51').Value = '  +2.57%  '

# Restore default (General) style on cells that needed a leading apostrophe
# to avoid Excel auto-converting numeric-looking text, so no quotePrefix/style
# drift is introduced relative to the original formatting.
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
